$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force every touched cell to remain plain text (matches source inlineStr cells),
# so numeric-looking strings like "0.549" or "39.01" are not coerced to floats.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '53.400.98'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.74%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.156.48'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.42%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '108.61'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +5.36%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.549'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.98%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +4.29%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.01'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +6.17%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0873'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.32%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.650.68'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.15%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.98%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.01'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.13%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +8.86%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.156.32'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.60%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '53.314.62'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +4.11%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.46%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0977'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.88'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '271.00'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.25'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.74%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.18'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '27.74'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.37'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.170'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.111'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.99%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.00'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +7.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '37.36'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +8.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0496'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +10.98%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.52%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '50.48'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.66'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +9.74%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.81'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +10.09%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +11.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.293'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.72%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '130.75'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +4.60%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.44%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.47'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.95%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.093.32'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.39'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.07%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0502'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +22.72%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.69'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +5.78%  '
